$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.293.07'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '2.244.62'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0812'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '2.336.85'
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("D15").Value = '2.588.14'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.832'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("D18").Value = '44.092.76'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '0.0₃0965'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.63'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("E29").Value = '  +2.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.88'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0793'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.87%  '
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("E38").Value = '  -7.84%  '
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.79%  '
$ws.Range("E42").Value = '  -1.67%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = '1.752.44'
$ws.Range("E44").Value = '  +2.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '99.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.58%  '
